$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 91; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = $val * 100
    }
}
